# PF-2242 - View potential partners as a final step of workflows (all)
# * Adding Taxonomy to the scoring
# * Giving two setups, one of which shows how multiple sub-scores works

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet (removes all existing rows/content/row-heights
# while keeping column definitions, styles table, etc. intact).
$ws.Rows("1:100").Delete()

# ---------------------------------------------------------------------
# Title / header block
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Partner Finder Matchmaker Algorithm Worksheet"
$ws.Cells.Item(2,1).Font.Bold = $true
$ws.Cells.Item(2,1).Font.Size = 20

$ws.Range("A3").Value = "v1.1"
$ws.Range("B3").Value = 44056
$ws.Range("B3").NumberFormat = "m/d/yy"
$ws.Range("C3").Value = 0.72916666666666663
$ws.Range("C3").NumberFormat = "h:mm AM/PM"

$ws.Range("A5").Value = "* Final Matchmaker Scores are 0.0 to 1.0, with 0 being the lowest quality match, and 1.0 being the best quality match"
$ws.Range("A6").Value = "* Final Matchmaker Scores are bidirectional, so a fitness score of 0.85 for Organization => Project is also a fitness score of 0.85 Project => Organization"
$ws.Range("A7").Value = "* There is a cutoff value in the system below which Final Matchmaker Scores won''t be displayed, MatchScoreDisplayCutoff. Right now this is 0.5. "
$ws.Range("A8").Value = "* Final Matchmaker Scores are composed of SubScores, which also range from 0.0 to 1.0. They are summed as follows:"

# ---------------------------------------------------------------------
# First example ("If/when we have multiple sub-scores:")
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "If/when we have multiple sub-scores:"
$ws.Cells.Item(11,1).Font.Bold = $true

$ws.Range("B12").Value = "SubScore Label"
$ws.Cells.Item(12,2).Font.Bold = $true
$ws.Range("C12").Value = "SubScore Score"
$ws.Cells.Item(12,3).Font.Bold = $true
$ws.Range("D12").Value = "Weighted SubScore"
$ws.Cells.Item(12,4).Font.Bold = $true

$ws.Range("B13").Value = "Taxonomy"
$ws.Range("C13").Formula = "=+`$B`$40"
$ws.Range("C13").NumberFormat = "0.00"
$ws.Range("D13").Formula = "=+C13*`$H`$15"
$ws.Range("D13").NumberFormat = "0.0000"
$ws.Range("G13").Value = "SubScore Count"
$ws.Cells.Item(13,7).Font.Bold = $true
$ws.Range("H13").Formula = "=COUNT(C13:C21)"

$ws.Range("B14").Value = "xxx"
$ws.Range("C14").Value = 0.1
$ws.Range("C14").NumberFormat = "0.00"
$ws.Range("D14").Formula = "=+C14*`$H`$15"
$ws.Range("D14").NumberFormat = "0.0000"
$ws.Range("G14").Value = "SubScore Max Value"
$ws.Cells.Item(14,7).Font.Bold = $true
$ws.Range("H14").Value = 1
$ws.Range("H14").NumberFormat = "0.0"

$ws.Range("B15").Value = "yyy"
$ws.Range("C15").Value = 0.2
$ws.Range("C15").NumberFormat = "0.00"
$ws.Range("D15").Formula = "=+C15*`$H`$15"
$ws.Range("D15").NumberFormat = "0.0000"
$ws.Range("G15").Value = "Per SubScore Weight"
$ws.Cells.Item(15,7).Font.Bold = $true
$ws.Range("H15").Formula = "=+H14/H13"
$ws.Range("H15").ClearFormats()

$ws.Range("B16").Value = "zzz"
$ws.Range("C16").Value = 0.3
$ws.Range("C16").NumberFormat = "0.00"
$ws.Range("D16").Formula = "=+C16*`$H`$15"
$ws.Range("D16").NumberFormat = "0.0000"

$ws.Range("C17").NumberFormat = "0.00"
$ws.Range("D17").NumberFormat = "0.0000"

$ws.Range("B18").Value = "Final Matchmaker Score"
$ws.Cells.Item(18,2).Font.Bold = $true
$ws.Cells.Item(18,2).Font.Size = 14
$ws.Cells.Item(18,3).Font.Bold = $true
$ws.Cells.Item(18,3).Font.Size = 14
$ws.Range("D18").Formula = "=SUM(D13:D16)"
$ws.Range("D18").NumberFormat = "0.0000"
$ws.Cells.Item(18,4).Font.Bold = $true
$ws.Cells.Item(18,4).Font.Size = 14
$ws.Rows(18).RowHeight = 18.75

# ---------------------------------------------------------------------
# Second example ("As implemented 8/13/2020 (only one subscore so far):")
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "As implemented 8/13/2020 (only one subscore so far):"
$ws.Cells.Item(20,1).Font.Bold = $true

$ws.Range("B22").Value = "SubScore Label"
$ws.Cells.Item(22,2).Font.Bold = $true
$ws.Range("C22").Value = "SubScore Score"
$ws.Cells.Item(22,3).Font.Bold = $true
$ws.Range("D22").Value = "Weighted SubScore"
$ws.Cells.Item(22,4).Font.Bold = $true

$ws.Range("B23").Value = "Taxonomy"
$ws.Range("C23").Formula = "=+`$B`$40"
$ws.Range("C23").NumberFormat = "0.00"
$ws.Range("D23").Formula = "=+C23*`$H`$25"
$ws.Range("D23").NumberFormat = "0.0000"
$ws.Range("G23").Value = "SubScore Count"
$ws.Cells.Item(23,7).Font.Bold = $true
$ws.Range("H23").Formula = "=COUNT(C23:C32)"

$ws.Range("C24").NumberFormat = "0.00"
$ws.Range("D24").NumberFormat = "0.0000"
$ws.Range("G24").Value = "SubScore Max Value"
$ws.Cells.Item(24,7).Font.Bold = $true
$ws.Range("H24").Value = 1
$ws.Range("H24").NumberFormat = "0.0"

$ws.Range("B25").Value = "Final Matchmaker Score"
$ws.Cells.Item(25,2).Font.Bold = $true
$ws.Cells.Item(25,2).Font.Size = 14
$ws.Cells.Item(25,3).Font.Bold = $true
$ws.Cells.Item(25,3).Font.Size = 14
$ws.Range("D25").Formula = "=SUM(D23)"
$ws.Range("D25").NumberFormat = "0.0000"
$ws.Cells.Item(25,4).Font.Bold = $true
$ws.Cells.Item(25,4).Font.Size = 14
$ws.Rows(25).RowHeight = 18.75
$ws.Range("G25").Value = "Per SubScore Weight"
$ws.Cells.Item(25,7).Font.Bold = $true
$ws.Range("H25").Formula = "=+H24/H23"
$ws.Range("H25").ClearFormats()

$ws.Range("G27").Font.Bold = $true
$ws.Range("G28").Font.Bold = $true

$ws.Range("A29").Value = "The per sub-score weighting is assumed to be even per sub-score above, but this is readily changed if we need."

# ---------------------------------------------------------------------
# Taxonomy SubScore detail / working area
# ---------------------------------------------------------------------
$ws.Range("A31").Value = "Taxonomy SubScore:"
$ws.Cells.Item(31,1).Font.Bold = $true
$ws.Cells.Item(31,1).Font.Size = 18
$ws.Rows(31).RowHeight = 23.25

$ws.Range("A32").Value = "(Explanation of Taxonomy SubScore could go here.)"
$ws.Cells.Item(32,1).Font.Italic = $true

$ws.Range("B33").Value = " "

$ws.Range("A34").Value = "Matches Leaf"
$ws.Range("B34").Value = $false

$ws.Range("A35").Value = "Matches Leaf on Branch"
$ws.Range("B35").Value = $true

$ws.Range("A36").Value = "Matches Leaf on Trunk"
$ws.Range("B36").Value = $false

$ws.Range("A38").Value = "Any Match"
$ws.Range("B38").Formula = "=OR(B34,B35,B36)"

$ws.Range("A40").Value = "Taxonomy SubScore"
$ws.Cells.Item(40,1).Font.Bold = $true
$ws.Range("B40").Formula = "=IF(B38, 1, 0)"
$ws.Range("B40").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# xxx / yyy SubScore placeholders
# ---------------------------------------------------------------------
$ws.Range("A42").Value = "xxx SubScore:"
$ws.Cells.Item(42,1).Font.Bold = $true
$ws.Cells.Item(42,1).Font.Size = 18
$ws.Rows(42).RowHeight = 23.25

$ws.Range("A43").Value = "(Explanation of xxx SubScore could go here.)"
$ws.Cells.Item(43,1).Font.Italic = $true

$ws.Range("A46").Value = "yyy SubScore:"
$ws.Cells.Item(46,1).Font.Bold = $true
$ws.Cells.Item(46,1).Font.Size = 18
$ws.Rows(46).RowHeight = 23.25

$ws.Range("A47").Value = "(Explanation of yyy SubScore could go here.)"
$ws.Cells.Item(47,1).Font.Italic = $true

# ---------------------------------------------------------------------
# Row heights for title / bullet rows
# ---------------------------------------------------------------------
$ws.Rows(2).RowHeight = 26.25

# ---------------------------------------------------------------------
# Column widths / sheet view
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 32.6

$ws.Range("C4").Select()
